$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# "Finish Find the Duplicate" - update the hours logged and notes for 2021-11-15 (row 70)
# Hours: 0.25 -> 1
$ws.Range("C70").Value = 1
# Notes + Weekly Total: "1 small problem" -> "2 small problems"
$ws.Range("D70").Value = "2 small problems"

# Move the active selection to D71 (next entry row)
$ws.Range("D71").Select()
